# ex9.1.9(Linear) - "expermits todos no convexos menos el 5to"
# Regenerate the non-convex follower restrictions (all but the 5th), the
# modified point, vec_bf, vec_BF and vec_alpha with the new generator output.

$wb = $excel.ActiveWorkbook

# Helper: force a value to be written as TEXT (matches the workbook's
# existing convention where every numeric-looking value in these tables is
# actually stored as a shared string, t="s"). We do this by writing a
# formula that evaluates to the literal text, then collapsing the formula
# down to its static value with a values-only paste - this avoids Excel's
# "looks like a number -> store as number" auto-conversion without leaving
# behind any NumberFormat/quotePrefix style residue.
function Set-TextValue($range, [string]$text) {
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
}

# ---------------------------------------------------------------------
# Restricciones_del_follower
# ---------------------------------------------------------------------
$wsF = $wb.Worksheets.Item("Restricciones_del_follower")

Set-TextValue $wsF.Range("A2") "-12.875620327718131 - x + 3.0295976275849865y"
Set-TextValue $wsF.Range("B2") "14.875620327718131"
Set-TextValue $wsF.Range("D2") "0.62"
Set-TextValue $wsF.Range("E2") "0"
Set-TextValue $wsF.Range("F2") "3.5999999999999996"

Set-TextValue $wsF.Range("A3") "-25.163290564403454 - 0.25x + 5.217640358618588y"
Set-TextValue $wsF.Range("B3") "23.163290564403454"
Set-TextValue $wsF.Range("D3") "0.96"
Set-TextValue $wsF.Range("E3") "0"
Set-TextValue $wsF.Range("F3") "6.2"

Set-TextValue $wsF.Range("A4") "-1.5668395382394484 + x - 0.0941383052358844y"
Set-TextValue $wsF.Range("B4") "-6.433160461760552"
Set-TextValue $wsF.Range("D4") "0.88"
Set-TextValue $wsF.Range("E4") "4.1"
Set-TextValue $wsF.Range("F4") "0.8"

Set-TextValue $wsF.Range("A5") "-2.85 + x"
Set-TextValue $wsF.Range("B5") "0.03000000000000025"
Set-TextValue $wsF.Range("D5") "0.29"
Set-TextValue $wsF.Range("E5") "6.4"
Set-TextValue $wsF.Range("F5") "0"

Set-TextValue $wsF.Range("A6") "-17.575845373234547 + 3.4503750758606797y"
Set-TextValue $wsF.Range("B6") "16.975845373234545"
Set-TextValue $wsF.Range("D6") "0.34"
Set-TextValue $wsF.Range("E6") "9.8"
Set-TextValue $wsF.Range("F6") "4.1"

$rngF = $wsF.Range("A2:F6")
$rngF.Copy()
$rngF.PasteSpecial(-4163)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Punto_modificado
# ---------------------------------------------------------------------
$wsP = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $wsP.Range("A2") "2.0300000000000002"
Set-TextValue $wsP.Range("B2") "4.92"
$rngP = $wsP.Range("A2:B2")
$rngP.Copy()
$rngP.PasteSpecial(-4163)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Vector_bf (lowercase name) - NOTE: both sheet NAMES ("Vector_bf" vs
# "Vector_BF") and PowerShell VARIABLE names are case-insensitive, so we
# must address these two worksheets by their (distinct) tab index rather
# than by name, or the lookups/vars collide with each other.
# ---------------------------------------------------------------------
$wsVecLower = $wb.Worksheets.Item(5)   # Vector_bf
Set-TextValue $wsVecLower.Range("A2") "-6.9775710905615895"
$rngVecLower = $wsVecLower.Range("A2")
$rngVecLower.Copy()
$rngVecLower.PasteSpecial(-4163)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Vector_BF (uppercase name)
# ---------------------------------------------------------------------
$wsVecUpper = $wb.Worksheets.Item(6)   # Vector_BF
Set-TextValue $wsVecUpper.Range("A2") "-11.5"
Set-TextValue $wsVecUpper.Range("A3") "-34.42770869196754"
$rngVecUpper = $wsVecUpper.Range("A2:A3")
$rngVecUpper.Copy()
$rngVecUpper.PasteSpecial(-4163)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Vector_Alpha (this one is a genuine number, not text)
# ---------------------------------------------------------------------
$wsAlpha = $wb.Worksheets.Item("Vector_Alpha")
$wsAlpha.Range("A2").Value = 1.1882766104717688
